# Update crypto price/volume figures per the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.595.06'
$ws.Range('E2').Value = '  +1.26%  '
$ws.Range('D3').Value = '3.021.47'
$ws.Range('E3').Value = '  +2.52%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '378.44'
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.89'
$ws.Range('E6').Value = '  +2.07%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.594'
$ws.Range('E9').Value = '  +2.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.53'
$ws.Range('E10').Value = '  +1.11%  '
$ws.Range('E11').Value = '  -0.36%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0860'
$ws.Range('E12').Value = '  +1.33%  '
$ws.Range('D13').Value = '3.496.08'
$ws.Range('E13').Value = '  +2.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.47'
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.73'
$ws.Range('E15').Value = '  +0.41%  '
$ws.Range('D16').Value = '3.020.46'
$ws.Range('E16').Value = '  +2.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.977'
$ws.Range('E17').Value = '  -1.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.67'
$ws.Range('E18').Value = '  -10.23%  '
$ws.Range('D19').Value = '51.600.56'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.03'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.47'
$ws.Range('E21').Value = '  +0.71%  '
$ws.Range('D22').Value = '0.0₃0960'
$ws.Range('E22').Value = '  +1.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.92'
$ws.Range('E23').Value = '  +0.66%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '269.00'
$ws.Range('E24').Value = '  +0.94%  '
$ws.Range('E25').Value = '  -3.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.27'
$ws.Range('E26').Value = '  +1.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.66'
$ws.Range('E27').Value = '  +8.58%  '
$ws.Range('E28').Value = '  +5.32%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '26.31'
$ws.Range('E29').Value = '  +2.71%  '
$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('E31').Value = '  +0.47%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.28'
$ws.Range('E32').Value = '  +2.04%  '
$ws.Range('B33').Value = 'VeChain'
$ws.Range('C33').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0457'
$ws.Range('E33').Value = '  +5.93%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '34.09'
$ws.Range('E34').Value = '  +1.89%  '
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.29'
$ws.Range('E38').Value = '  +6.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.289'
$ws.Range('E39').Value = '  +10.92%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.11'
$ws.Range('E40').Value = '  +3.42%  '
$ws.Range('E41').Value = '  +4.09%  '
$ws.Range('E42').Value = '  +2.84%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '127.41'
$ws.Range('E43').Value = '  +5.97%  '
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('E45').Value = '  +8.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '21.65'
$ws.Range('E46').Value = '  +1.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.08'
$ws.Range('E47').Value = '  +3.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.39'
$ws.Range('E48').Value = '  +2.70%  '
$ws.Range('D49').Value = '2.029.54'
$ws.Range('E49').Value = '  +0.93%  '
$ws.Range('D50').Value = '3.320.00'
$ws.Range('E50').Value = '  +2.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0320'
$ws.Range('E51').Value = '  +2.12%  '
